$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel will not auto-convert to numbers)
$ws.Range("D2").Value = '42.227.57'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.246.34'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("E7").Value = '  -5.63%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -3.67%  '
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("E11").Value = '  -2.54%  '
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("E14").Value = '  -4.79%  '
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").Value = '2.258.23'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '42.118.11'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '0.0₃0995'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("E23").Value = '  +36.90%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("E28").Value = '  +3.57%  '
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("E31").Value = '  -4.30%  '
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E34").Value = '  -0.86%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E35").Value = '  +10.02%  '
$ws.Range("E36").Value = '  -2.74%  '
$ws.Range("E37").Value = '  +3.56%  '
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("E39").Value = '  -3.95%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("E42").Value = '  -2.53%  '
$ws.Range("E43").Value = '  -6.77%  '
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("E45").Value = '  -2.90%  '
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("E50").Value = '  -8.69%  '
$ws.Range("E51").Value = '  -0.23%  '

# Numeric-looking text updates: force Text format so Excel keeps them as strings
# (e.g. "74.50" must stay "74.50", not become the number 74.5),
# then reset the cell style back to Normal so no stray style/number-format is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.854'
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0819'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.125'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0315'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '13.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '62.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.205'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.21'
$ws.Range("D50").Style = "Normal"
